$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New test-case rows appended to the check list (rows 19-21)
$newRows = @(
    @{ Row = 19; Num = 1.18; Text = "Проверка отображения полей с большим количеством символов"; Height = 34 },
    @{ Row = 20; Num = 1.19; Text = "Отображение пустой таблицы"; Height = 17 },
    @{ Row = 21; Num = 1.21; Text = "отображение большого количества записей"; Height = 17 }
)

foreach ($r in $newRows) {
    $rowIndex = $r.Row

    # Copy formatting from the row above (row 18) so the new rows match
    # the rest of the check list exactly.
    $ws.Range("A18:C18").Copy() | Out-Null
    $ws.Range("A$rowIndex`:C$rowIndex").PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($rowIndex, 1).Value = $r.Num
    $ws.Cells.Item($rowIndex, 2).Value = $r.Text
    $ws.Cells.Item($rowIndex, 3).Value = "passed"

    $ws.Rows.Item($rowIndex).RowHeight = $r.Height
}

$excel.CutCopyMode = 0

$ws.Range("E20").Select()
